{"js": "// Rewrite the \"KEY ACHIEVEMENTS AND IMPACT\" > \"Impact\" bullet list so that\n// the six job-duty-style bullets become four impact-focused accomplishment\n// statements, per the commit's \"Fix Key Achievements to use proper\n// accomplishment statements\" change.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the six existing bullets under KEY ACHIEVEMENTS AND IMPACT by their\n// exact current text, so the edit is anchored to content rather than a\n// fragile fixed index.\nconst oldBullets = [\n  \"\\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%\",\n  \"\\u2022 Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations\",\n  \"\\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis\",\n  \"\\u2022 Developed longitudinal data analysis methods using geospatial techniques that improved segmentation accuracy by 34% and survey incidence rates by 28%, reducing polling costs while increasing response quality\",\n  \"\\u2022 Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets\",\n  \"\\u2022 Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy\"\n];\n\nconst newBullets = [\n  \"\\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\",\n  \"\\u2022 $4.7M savings enabled nonprofit access\",\n  \"\\u2022 Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\",\n  \"\\u2022 178% accuracy improvement in racial classification algorithms\"\n];\n\n// Find the first paragraph index where all six old bullets appear\n// consecutively (this section is the only place they appear together, right\n// after the \"Impact\" sub-heading under KEY ACHIEVEMENTS AND IMPACT).\nconst items = paragraphs.items;\nlet startIdx = -1;\nfor (let i = 0; i + oldBullets.length <= items.length; i++) {\n  let allMatch = true;\n  for (let j = 0; j < oldBullets.length; j++) {\n    if (items[i + j].text !== oldBullets[j]) {\n      allMatch = false;\n      break;\n    }\n  }\n  if (allMatch) {\n    startIdx = i;\n    break;\n  }\n}\n\nif (startIdx === -1) {\n  throw new Error(\"Could not locate the six Key Achievements bullets to replace.\");\n}\n\n// Overwrite the first four matched paragraphs with the new bullet text.\nfor (let j = 0; j < newBullets.length; j++) {\n  items[startIdx + j].getRange().insertText(newBullets[j], Word.InsertLocation.replace);\n}\n\n// Remove the trailing two paragraphs (five and six) that no longer have a\n// corresponding new bullet.\nfor (let j = newBullets.length; j < oldBullets.length; j++) {\n  items[startIdx + j].delete();\n}\n\nawait context.sync();\n", "ps1": "# Rewrite the \"KEY ACHIEVEMENTS AND IMPACT\" > \"Impact\" bullet list so that\n# the six job-duty-style bullets become four impact-focused accomplishment\n# statements, per the commit's \"Fix Key Achievements to use proper\n# accomplishment statements\" change.\n\n$d = $word.ActiveDocument\n\n$bullet = [char]0x2022\n\n$oldBullets = @(\n  \"$bullet Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%\",\n  \"$bullet Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations\",\n  \"$bullet Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis\",\n  \"$bullet Developed longitudinal data analysis methods using geospatial techniques that improved segmentation accuracy by 34% and survey incidence rates by 28%, reducing polling costs while increasing response quality\",\n  \"$bullet Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets\",\n  \"$bullet Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy\"\n)\n\n$newBullets = @(\n  \"$bullet Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\",\n  \"$bullet `$4.7M savings enabled nonprofit access\",\n  \"$bullet Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\",\n  \"$bullet 178% accuracy improvement in racial classification algorithms\"\n)\n\n# Find the first paragraph index (1-based, like Word's Paragraphs collection)\n# where all six old bullets appear consecutively - this section is the only\n# place they appear together, right after the \"Impact\" sub-heading under\n# KEY ACHIEVEMENTS AND IMPACT.\n$count = $d.Paragraphs.Count\n$startIdx = -1\nfor ($i = 1; $i -le ($count - $oldBullets.Count + 1); $i++) {\n  $allMatch = $true\n  for ($j = 0; $j -lt $oldBullets.Count; $j++) {\n    $paraText = $d.Paragraphs.Item($i + $j).Range.Text\n    $paraText = $paraText.TrimEnd([char]13, [char]7)\n    if ($paraText -ne $oldBullets[$j]) {\n      $allMatch = $false\n      break\n    }\n  }\n  if ($allMatch) {\n    $startIdx = $i\n    break\n  }\n}\n\nif ($startIdx -eq -1) {\n  throw \"Could not locate the six Key Achievements bullets to replace.\"\n}\n\n# Overwrite the first four matched paragraphs with the new bullet text\n# (exclude the trailing paragraph mark from the replaced range so no extra\n# paragraph break gets inserted).\nfor ($j = 0; $j -lt $newBullets.Count; $j++) {\n  $r = $d.Paragraphs.Item($startIdx + $j).Range\n  $r.MoveEnd(1, -1) | Out-Null\n  $r.Text = $newBullets[$j]\n}\n\n# Remove the trailing two paragraphs (five and six) that no longer have a\n# corresponding new bullet. Delete from the back so earlier indices stay\n# valid.\nfor ($j = $oldBullets.Count - 1; $j -ge $newBullets.Count; $j--) {\n  $d.Paragraphs.Item($startIdx + $j).Range.Delete()\n}\n"}
